$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts old rows 33-102 down to 34-103,
# carrying their formatting/values with them (matches the target diff which
# shows the entire former row N re-appearing as row N+1 for N = 33..102).
$ws.Rows(33).Insert()

# Populate the newly inserted row 33 with the new weekly price record.
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 44607
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 100112001
$ws.Cells.Item(33, 7).Value = "Berenjena"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 7000
$ws.Cells.Item(33, 12).Value = 7000
$ws.Cells.Item(33, 13).Value = 7000
$ws.Cells.Item(33, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(33, 15).Value = "Región del Maule"
$ws.Cells.Item(33, 16).Value = 140
$ws.Cells.Item(33, 17).Value = 50
$ws.Cells.Item(33, 18).Value = "Hortaliza"
